$d = $word.ActiveDocument

function Find-ParagraphIndex($pattern) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -match $pattern) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# Edit 1: the "Performance Requirements" narrative paragraph had its final
# sentence split across a stray "_GoBack" bookmark into two runs; merge them
# back into one run and drop the bookmark.
# ---------------------------------------------------------------------------
$perfIdx = Find-ParagraphIndex("^Performance requirements define sustainable")
if ($perfIdx -eq -1) {
    throw "Could not find the Performance Requirements narrative paragraph"
}
$perfRange = $d.Paragraphs.Item($perfIdx).Range.Duplicate

$perfXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:ind w:left="-5" w:firstLine="487"/><w:jc w:val="both"/><w:rPr><w:i w:val="0"/></w:rPr></w:pPr>
<w:r><w:rPr><w:i w:val="0"/></w:rPr><w:t>Performance</w:t></w:r>
<w:r><w:rPr><w:i w:val="0"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:rPr><w:i w:val="0"/></w:rPr><w:t xml:space="preserve">requirements define sustainable response times for system functionality. </w:t></w:r>
<w:r><w:rPr><w:i w:val="0"/></w:rPr><w:t xml:space="preserve">Despite the fact that the system is created suiting for the least system performances, the performance of the system will highly depend on the performance of the hardware and software components of the installing computer. At the point when consider about the timing relationships of the system, the load time for </w:t></w:r>
<w:r><w:rPr><w:i w:val="0"/></w:rPr><w:t xml:space="preserve">user </w:t></w:r>
<w:r><w:rPr><w:i w:val="0"/></w:rPr><w:t>interface</w:t></w:r>
<w:r><w:rPr><w:i w:val="0"/></w:rPr><w:t xml:space="preserve"> screens might take no longer</w:t></w:r>
<w:r><w:rPr><w:i w:val="0"/></w:rPr><w:t xml:space="preserve"> than two</w:t></w:r>
<w:r><w:rPr><w:i w:val="0"/></w:rPr><w:t xml:space="preserve"> seconds. It makes fast access to system functions. The log in details shall be verified within five seconds causes&#8217; efficiency of the system. Returning query results within five seconds makes search function more accurate.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$perfRange.InsertXML($perfXml)

# ---------------------------------------------------------------------------
# Edit 2: the "Safety Requirements" section's placeholder paragraph
# ("<Specify those requirements ...>") is replaced with real content
# describing user levels / login / backups, and it picks up the relocated
# "_GoBack" bookmark plus justified / first-line-indented formatting.
# ---------------------------------------------------------------------------
$placeholderIdx = Find-ParagraphIndex("^<Specify those requirements that are concerned with possible loss")
if ($placeholderIdx -eq -1) {
    throw "Could not find the Safety Requirements placeholder paragraph"
}
$safetyRange = $d.Paragraphs.Item($placeholderIdx).Range.Duplicate

$safetyXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:ind w:left="-5" w:firstLine="487"/>
<w:jc w:val="both"/>
<w:rPr><w:i w:val="0"/></w:rPr>
</w:pPr>
<w:bookmarkStart w:id="1001" w:name="_GoBack"/>
<w:bookmarkEnd w:id="1001"/>
<w:r><w:rPr><w:i w:val="0"/></w:rPr><w:t xml:space="preserve">There are several user levels in resort reservation system, access to the different subsystems will be ensured by a user log in screen that requires a username and password. This </w:t></w:r>
<w:r><w:rPr><w:i w:val="0"/></w:rPr><w:t>gives different perspectives</w:t></w:r>
<w:r><w:rPr><w:i w:val="0"/></w:rPr><w:t xml:space="preserve"> and</w:t></w:r>
<w:r><w:rPr><w:i w:val="0"/></w:rPr><w:t xml:space="preserve"> accessible functions of user levels through the system. Maintaining backups ensure the system database security. System can be restoring in any case of emergency.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$safetyRange.InsertXML($safetyXml)

Write-Output "Applied Safety Requirements edits"
